# Update the "TestData7"-related duration figures (F/G columns) for rows
# 4, 6, 7, 8, 9, 10 on Sheet1: the F column value moves from "16" to "17"
# and the G column value moves from "14" to "15". These cells are stored
# as text (quote-prefixed numerics), so keep them as text here too.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "17"
    $ws.Cells.Item($r, 7).Value = "15"
}

# Move the active selection from L10 to H12.
$ws.Range("H12").Select()
